# edit.ps1 - applies the two changes described by the diff:
#   1. After the "Justin Bilao" paragraph (title-page author list), insert
#      a blank bold paragraph followed by a bold "Test upload to main"
#      paragraph.
#   2. Merge the two runs that make up "User-Defined Period Analysis:" +
#      " " into a single run "User-Defined Period Analysis: ".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: insert the two new paragraphs after "Justin Bilao".
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Justin Bilao", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)

    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

    $xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="1C0EFEAB" w14:textId="726AF29E" w:rsidR="00BB6FD3" w:rsidRPr="00BB6FD3" w:rsidRDefault="00BB6FD3" w:rsidP="00926CFD">' +
             '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
             '<w:r w:rsidRPr="00BB6FD3"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Justin Bilao</w:t></w:r>' +
           '</w:p>' +
           '<w:p ' + $wNs + '>' +
             '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
           '</w:p>' +
           '<w:p ' + $wNs + '>' +
             '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
             '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Test upload to main</w:t></w:r>' +
           '</w:p>'

    $rng.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------
# Change 2: collapse "User-Defined Period Analysis:" + " " (two runs)
# into a single run "User-Defined Period Analysis: ".
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("User-Defined Period Analysis: ", $true, $false, $false, $false, $false, $true, 1, $false, "User-Defined Period Analysis: ", 2)
